$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 935.2857
$ws.Range("I6").Value = 224.58333
$ws.Range("K6").Value = 673.74999
$ws.Range("M6").Value = -561.74999

# Row 8
$ws.Range("H8").Value = 60.5
$ws.Range("I8").Value = 60.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 181.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -42.5
$ws.Range("N8").ClearContents()

# Row 38
$ws.Range("H38").Value = 371
$ws.Range("I38").Value = 44.285713
$ws.Range("K38").Value = 132.857139
$ws.Range("M38").Value = 239.142861

# Row 39
$ws.Range("H39").Value = 459.42856
$ws.Range("I39").Value = 603.2
$ws.Range("J39").Value = 100
$ws.Range("K39").Value = 1809.6
$ws.Range("L39").Value = 300
$ws.Range("M39").Value = -1513.6
$ws.Range("N39").Value = -892

# Row 55
$ws.Range("H55").Value = 181
$ws.Range("I55").Value = 114.57143
$ws.Range("J55").Value = 297.25
$ws.Range("K55").Value = 114.57143
$ws.Range("L55").Value = 297.25
$ws.Range("M55").Value = 99.42856999999999
$ws.Range("N55").Value = -725.25

# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 98
$ws.Range("H98").Value = 1598.7778
$ws.Range("I98").Value = 647.5
$ws.Range("J98").Value = 2359.8
$ws.Range("K98").Value = 647.5
$ws.Range("L98").Value = 2359.8
$ws.Range("M98").Value = 850.5
$ws.Range("N98").Value = -5355.8

# Row 112
$ws.Range("H112").Value = 3103.6875
$ws.Range("J112").Value = 3217.4
$ws.Range("L112").Value = 9652.200000000001
$ws.Range("N112").Value = -11868.2

# Row 116
$ws.Range("H116").Value = 5206.7144
$ws.Range("I116").Value = 5139.6
$ws.Range("J116").Value = 5374.5
$ws.Range("K116").Value = 5139.6
$ws.Range("L116").Value = 5374.5
$ws.Range("M116").Value = -1697.6
$ws.Range("N116").Value = -12258.5

# Row 122
$ws.Range("H122").Value = 1598.7778
$ws.Range("I122").Value = 647.5
$ws.Range("J122").Value = 2359.8
$ws.Range("K122").Value = 1942.5
$ws.Range("L122").Value = 7079.400000000001
$ws.Range("M122").Value = 507.5
$ws.Range("N122").Value = -11979.4

# Row 132
$ws.Range("H132").Value = 2209.2068
$ws.Range("I132").Value = 1310.3846
$ws.Range("K132").Value = 3931.1538
$ws.Range("M132").Value = -1401.1538

# Row 137
$ws.Range("H137").Value = 600
$ws.Range("I137").Value = 600
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 1800
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = 750
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 8600.666999999999
$ws.Range("J37").Value = 20000
$ws.Range("L37").Value = 20000
$ws.Range("N37").Value = -20546

# Row 45
$ws.Range("H45").Value = 2976.2
$ws.Range("I45").Value = 1160.75
$ws.Range("K45").Value = 1160.75
$ws.Range("M45").Value = -783.75

# Row 61
$ws.Range("H61").Value = 10102.6
$ws.Range("J61").Value = 19628.25
$ws.Range("L61").Value = 19628.25
$ws.Range("N61").Value = -20052.25

# Row 74
$ws.Range("H74").Value = 3499.5
$ws.Range("I74").Value = 3499.5
$ws.Range("K74").Value = 3499.5
$ws.Range("M74").Value = -2625.5

# Row 77
$ws.Range("H77").Value = 3499.5
$ws.Range("I77").Value = 3499.5
$ws.Range("K77").Value = 17497.5
$ws.Range("M77").Value = -13129.5

# Row 136
$ws.Range("H136").Value = 10102.6
$ws.Range("J136").Value = 19628.25
$ws.Range("L136").Value = 58884.75
$ws.Range("N136").Value = -63984.75

$ws = $wb.Worksheets.Item("BSM")
# Row 5
$ws.Range("H5").Value = 412.5
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 325
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 325
$ws.Range("M5").Value = -387
$ws.Range("N5").Value = -551

# Row 20
$ws.Range("H20").Value = 1099.5
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 1099.5
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 1099.5
$ws.Range("N20").Value = -1593.5
$ws.Range("M20").ClearContents()

# Row 22
$ws.Range("H22").Value = 661.3333
$ws.Range("I22").Value = 661.3333
$ws.Range("K22").Value = 661.3333
$ws.Range("M22").Value = -488.3333

# Row 86
$ws.Range("H86").Value = 1665.2222
$ws.Range("I86").Value = 1926.8572
$ws.Range("K86").Value = 1926.8572
$ws.Range("M86").Value = -803.8571999999999

# Row 89
$ws.Range("H89").Value = 1665.2222
$ws.Range("I89").Value = 1926.8572
$ws.Range("K89").Value = 9634.286
$ws.Range("M89").Value = -4018.286

# Row 107
$ws.Range("H107").Value = 699.13635
$ws.Range("I107").Value = 699.0952
$ws.Range("K107").Value = 699.0952
$ws.Range("M107").Value = 1220.9048

$ws = $wb.Worksheets.Item("CRP")
# Row 12
$ws.Range("H12").Value = 25301.125
$ws.Range("J12").Value = 50000
$ws.Range("L12").Value = 50000
$ws.Range("N12").Value = -50340

# Row 69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

# Row 72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

# Row 107
$ws.Range("H107").Value = 2130.5454
$ws.Range("I107").Value = 1710
$ws.Range("J107").Value = 2738
$ws.Range("K107").Value = 1710
$ws.Range("L107").Value = 2738
$ws.Range("M107").Value = 210
$ws.Range("N107").Value = -6578

# Row 134
$ws.Range("H134").Value = 1904.7778
$ws.Range("I134").Value = 1886.4
$ws.Range("K134").Value = 5659.200000000001
$ws.Range("M134").Value = -3124.200000000001

$ws = $wb.Worksheets.Item("CUL")
# Row 41
$ws.Range("H41").Value = 525
$ws.Range("I41").Value = 300
$ws.Range("K41").Value = 900
$ws.Range("M41").Value = -562

# Row 60
$ws.Range("H60").Value = 3449.75
$ws.Range("I60").Value = 3599.6667
$ws.Range("J60").Value = 3000
$ws.Range("K60").Value = 10799.0001
$ws.Range("L60").Value = 9000
$ws.Range("M60").Value = -10548.0001
$ws.Range("N60").Value = -9502

# Row 76
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M76").ClearContents()

# Row 79
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("M79").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 46
$ws.Range("H46").Value = 7937.375
$ws.Range("J46").Value = 12000
$ws.Range("L46").Value = 12000
$ws.Range("N46").Value = -12312

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1577.8948
$ws.Range("I96").Value = 1508.6428
$ws.Range("J96").Value = 1771.8
$ws.Range("K96").Value = 1508.6428
$ws.Range("L96").Value = 1771.8
$ws.Range("M96").Value = -135.6428000000001
$ws.Range("N96").Value = -4517.8

# Row 107
$ws.Range("H107").Value = 255.77777
$ws.Range("I107").Value = 255.77777
$ws.Range("K107").Value = 767.33331
$ws.Range("M107").Value = 1152.66669

# Row 127
$ws.Range("H127").Value = 25000
$ws.Range("J127").Value = 25000
$ws.Range("L127").Value = 25000
$ws.Range("N127").Value = -34920

# Row 136
$ws.Range("H136").Value = 896.5333000000001
$ws.Range("I136").Value = 957.9231
$ws.Range("K136").Value = 2873.7693
$ws.Range("M136").Value = -323.7692999999999
